# Add a "Time Left" column (K) with auction time-remaining data,
# and refresh a few of the eBay URLs / price figures that changed
# on re-scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1, matching the style of the other header cells (B1:J1)
$ws.Range("K1").Value = "Time Left"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Row 2
$ws.Range("C2").Value = "https://www.ebay.com/itm/334664407065?hash=item4deb8d9819:g:FKIAAOSwnVBjmHe-&amdata=enc%3AAQAHAAAAoCG5TsTwv3NrcsAMph23YpyOzpjg%2FazSUpiMv3ki3a%2F8FVgE5ZK4lZIiQGyAC%2FdZ4xGQCBq3%2FmamHjwsFXxsfaUPLZ1sFmtJUXDunazGVqfFrr53b7Olk3lu9gylYT9DJgSxOvup5%2BkvNA859qS20gbDjdMJ6U8RU9SnjTuGxZi7G4vsWVnfwK2CQw3yM%2FnNxux1K%2Fe6QjyIyJEepHrwvJk%3D%7Ctkp%3ABk9SR8C8t7WjYQ"
$ws.Range("K2").Value = "Not specified"

# Row 3
$ws.Range("C3").Value = "https://www.ebay.com/itm/385163444573?epid=23050214728&hash=item59ad87c15d:g:oFAAAOSw5yxjRIsG&amdata=enc%3AAQAHAAAAoHHTmOg%2BoR3LtRCMXYQ8lAEEmjWes9DhmmiUWq2HMpjVVk%2FNDA2W4sdxSIQ5obiJZQz00dOOV%2FK5RwnOY3%2B7GF1HujRXbDUY7TxyD3o%2Fh4cwR6khxppUTdNw74CKaC%2FPwZHdFpdhHZmEnu29pPOZ9yGbKa501U6N703auECAI5PIi6s5BTezv1fm3sviBmhA0FvuvbQm0K26sjHtryL3Fmk%3D%7Ctkp%3ABk9SR8K8t7WjYQ"
$ws.Range("K3").Value = "Not specified"

# Row 4
$ws.Range("C4").Value = "https://www.ebay.com/itm/364076525731?hash=item54c4a6e8a3:g:JX4AAOSw~CJjlfzi&amdata=enc%3AAQAHAAAA4NEFSvBjjHf1HqZEVcXoZ9aEHAWu9bZtDogcLQ1aECJrJ4bqtAM2GcxfCgyJWBDbzJoIrzPMG0Cg7y%2BGOVailhMwXU4A7A09ESoKAM2JBOuXId253o%2F8uJGTEwKxMubxgfIgbB4jdrxHgk4q3WP7JKDvqFgkqWEHx0YOfpund8EmW30IjJIoBYFcvZx2%2B05boDpcuS%2BeVKLgZsfaWoVIV3gJubhgvE%2Fa8RjnHdJuEZhdmNu6sd9mMVAM4Q8w3p6hTzQCIkbucd4yHGzo8%2B%2FjPKdw0sVBzXJMx8oBJbkNGNQR%7Ctkp%3ABk9SR8C8t7WjYQ"
$ws.Range("E4").Value = 450
$ws.Range("H4").Value = 465
$ws.Range("K4").Value = "2d 1h restantes"

# Row 5
$ws.Range("C5").Value = "https://www.ebay.com/itm/155308020238?hash=item242914520e:g:HecAAOSw6CBjmLE8&amdata=enc%3AAQAHAAAAkINXBQ6TDT7aR4w2cDXsHla8%2F9Ul%2FpYukvcoYW7Eov4CrBmSjpfukU7mosoB2HqrTVLki2AttJRfLQmk4wQb%2BfDKJ4WiKr46Fc3erFfTH0FAEZJ%2BjBTl7ycxS3IXvsu0bHOUmPBmTmS%2FRBlcU2zxbwJkr9BFfARP0WC8nRTKDdGxCcRDVTgaLZ0Id2mmtsRjrA%3D%7Ctkp%3ABk9SR8C8t7WjYQ"
$ws.Range("E5").Value = 506
$ws.Range("H5").Value = 521
$ws.Range("K5").Value = "2d 2h restantes"

# Row 6
$ws.Range("C6").Value = "https://www.ebay.com/itm/175531638048?hash=item28de803d20:g:1rwAAOSwGZ5jjxIj&amdata=enc%3AAQAHAAAAkDUDHPwW2TkeJpPu0PLgguIKqM6OK8we5QyqR8AY203rTSvKqOyw6cxUU5TxcD25ytKRn5wIIAvNEsx%2Fua6UFblJPEmstQdP5a8RMKu9eMouJp0Chy9jf8qtXa86zTXO%2BHdjAqNn%2FpfUBS6jA1OCgMt57G%2F5BBQ2GMK9cAJjrIEA%2F0fLGXYewAECtRQGzVCWzg%3D%7Ctkp%3ABk9SR8K8t7WjYQ"
$ws.Range("K6").Value = "3d 20h restantes"
